# Scheduled market-data refresh: updates currentAveragePrice(NQ/HQ) and
# derived Leve cost/profit columns (H:N) on a handful of rows across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3981458
$ws.Range("I62").Value = 6546619
$ws.Range("J62").Value = 17118.182
$ws.Range("K62").Value = 6546619
$ws.Range("L62").Value = 17118.182
$ws.Range("M62").Value = -6545995
$ws.Range("N62").Value = -18366.182

$ws.Range("H65").Value = 3981458
$ws.Range("I65").Value = 6546619
$ws.Range("J65").Value = 17118.182
$ws.Range("K65").Value = 32733095
$ws.Range("L65").Value = 85590.91
$ws.Range("M65").Value = -32729975
$ws.Range("N65").Value = -91830.91

$ws.Range("H70").Value = 1961.15
$ws.Range("I70").Value = 1100
$ws.Range("J70").Value = 2248.2
$ws.Range("K70").Value = 3300
$ws.Range("L70").Value = 6744.599999999999
$ws.Range("M70").Value = -3030
$ws.Range("N70").Value = -7284.599999999999

$ws.Range("H73").Value = 1961.15
$ws.Range("I73").Value = 1100
$ws.Range("J73").Value = 2248.2
$ws.Range("K73").Value = 3300
$ws.Range("L73").Value = 6744.599999999999
$ws.Range("M73").Value = -2364
$ws.Range("N73").Value = -8616.599999999999

$ws.Range("H127").Value = 689
$ws.Range("I127").Value = 481.66666
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 1444.99998
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 3515.00002
$ws.Range("N127").Value = -12920

$ws.Range("H138").Value = 3907199
$ws.Range("I138").Value = 1138423.4
$ws.Range("J138").Value = 5850199.5
$ws.Range("K138").Value = 3415270.2
$ws.Range("L138").Value = 17550598.5
$ws.Range("M138").Value = -3410130.2
$ws.Range("N138").Value = -17560878.5

$ws.Range("H141").Value = 1910.4271
$ws.Range("I141").Value = 1089.9625
$ws.Range("J141").Value = 6012.75
$ws.Range("K141").Value = 3269.8875
$ws.Range("L141").Value = 18038.25
$ws.Range("M141").Value = 1910.1125
$ws.Range("N141").Value = -28398.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1085.2222
$ws.Range("I45").Value = 845.875
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 845.875
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -468.875
$ws.Range("N45").Value = -3754

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("N81").Value = 0

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("N84").Value = 0

$ws.Range("H122").Value = 16943.428
$ws.Range("I122").Value = 19117.334
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 57352.00199999999
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -54902.00199999999
$ws.Range("N122").Value = -16600

$ws.Range("H132").Value = 1651.6833
$ws.Range("I132").Value = 1245.3654
$ws.Range("J132").Value = 4292.75
$ws.Range("K132").Value = 3736.0962
$ws.Range("L132").Value = 12878.25
$ws.Range("M132").Value = -1206.0962
$ws.Range("N132").Value = -17938.25

$ws.Range("H141").Value = 63476.332
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 63476.332
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 63476.332
$ws.Range("N141").Value = -73836.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19610110
$ws.Range("I134").Value = 22729014
$ws.Range("J134").Value = 5573.143
$ws.Range("K134").Value = 68187042
$ws.Range("L134").Value = 16719.429
$ws.Range("M134").Value = -68184507
$ws.Range("N134").Value = -21789.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1481.5625
$ws.Range("I5").Value = 379.1
$ws.Range("J5").Value = 3319
$ws.Range("K5").Value = 379.1
$ws.Range("L5").Value = 3319
$ws.Range("M5").Value = -267.1
$ws.Range("N5").Value = -3543

$ws.Range("H31").Value = 1959.0238
$ws.Range("I31").Value = 1197.2667
$ws.Range("J31").Value = 3863.4167
$ws.Range("K31").Value = 1197.2667
$ws.Range("L31").Value = 3863.4167
$ws.Range("M31").Value = -902.2666999999999
$ws.Range("N31").Value = -4453.4167

$ws.Range("H34").Value = 1959.0238
$ws.Range("I34").Value = 1197.2667
$ws.Range("J34").Value = 3863.4167
$ws.Range("K34").Value = 1197.2667
$ws.Range("L34").Value = 3863.4167
$ws.Range("M34").Value = -995.2666999999999
$ws.Range("N34").Value = -4267.4167

$ws.Range("H58").Value = 1025.5172
$ws.Range("I58").Value = 706.54
$ws.Range("J58").Value = 3019.125
$ws.Range("K58").Value = 706.54
$ws.Range("L58").Value = 3019.125
$ws.Range("M58").Value = -503.54
$ws.Range("N58").Value = -3425.125

$ws.Range("H132").Value = 1716.8226
$ws.Range("I132").Value = 1368.4464
$ws.Range("J132").Value = 4968.3335
$ws.Range("K132").Value = 4105.3392
$ws.Range("L132").Value = 14905.0005
$ws.Range("M132").Value = -1575.3392
$ws.Range("N132").Value = -19965.0005

$ws.Range("H136").Value = 1025.5172
$ws.Range("I136").Value = 706.54
$ws.Range("J136").Value = 3019.125
$ws.Range("K136").Value = 2119.62
$ws.Range("L136").Value = 9057.375
$ws.Range("M136").Value = 430.3800000000001
$ws.Range("N136").Value = -14157.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 900
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").Value = 2700
$ws.Range("N16").Value = -3046

$ws.Range("H113").Value = 14706447
$ws.Range("I113").Value = 607.5
$ws.Range("J113").Value = 27778304
$ws.Range("K113").Value = 1822.5
$ws.Range("L113").Value = 83334912
$ws.Range("M113").Value = 347.5
$ws.Range("N113").Value = -83339252

$ws.Range("H128").Value = 404249.25
$ws.Range("I128").Value = 404249.25
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 1212747.75
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -1207767.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 34979
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 34979
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 34979
$ws.Range("N39").Value = -36043

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1833.3334
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2250
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 2250
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -2626

$ws.Range("H115").Value = 26792.5
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 26792.5
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 26792.5
$ws.Range("N115").Value = -29142.5

$ws.Range("H122").Value = 3671.7144
$ws.Range("I122").Value = 3202
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 9606
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -7156
$ws.Range("N122").Value = -16150

$ws.Range("H132").Value = 5243.4224
$ws.Range("I132").Value = 5220.926
$ws.Range("J132").Value = 5277.1665
$ws.Range("K132").Value = 15662.778
$ws.Range("L132").Value = 15831.4995
$ws.Range("M132").Value = -13132.778
$ws.Range("N132").Value = -20891.4995

$ws.Range("H136").Value = 3012.951
$ws.Range("I136").Value = 1513.66
$ws.Range("J136").Value = 9827.909
$ws.Range("K136").Value = 4540.98
$ws.Range("L136").Value = 29483.727
$ws.Range("M136").Value = -1990.98
$ws.Range("N136").Value = -34583.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 11500
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 11500
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 11500
$ws.Range("N29").Value = -12080

$ws.Range("H110").Value = 31011
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 31011
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 31011
$ws.Range("N110").Value = -39191

$ws.Range("H132").Value = 10002532
$ws.Range("I132").Value = 14708162
$ws.Range("J132").Value = 3069.875
$ws.Range("K132").Value = 44124486
$ws.Range("L132").Value = 9209.625
$ws.Range("M132").Value = -44121956
$ws.Range("N132").Value = -14269.625
